$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stage copies of the 8 pre-existing cell formats (old cellXfs indices
#    4..11) into a scratch row (100) far outside the used range, so we can
#    re-apply them (in permuted positions) without clobbering a source cell
#    before it has been read. xlPasteFormats = -4122.
#    A100<-old style 4 (A27), B100<-5 (B27), C100<-6 (D27), D100<-7 (A28),
#    E100<-8 (B28), F100<-9 (D28), G100<-10 (A8), H100<-11 (A1)
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$stageSources = @("A27", "B27", "D27", "A28", "B28", "D28", "A8", "A1")
$stageTargets = @("A100", "B100", "C100", "D100", "E100", "F100", "G100", "H100")
for ($i = 0; $i -lt $stageSources.Length; $i++) {
    $ws.Range($stageSources[$i]).Copy()
    $ws.Range($stageTargets[$i]).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# 2) Re-apply the rotated formats to their destination ranges.
#    (the style previously used by the "A1/A10/A19/A29/A39" date-header
#    band now belongs to the "A8/A17" summary band and vice-versa, and the
#    A27/B27/.../D28 bordered-label block rotates by one position)
# ---------------------------------------------------------------------------
$applyMap = @{
    "H100" = @("A1:D1", "A10:D10", "A19:D19", "A29:D29", "A39:D39")
    "G100" = @("A8:D9", "A17:D18")
    "A100" = @("A27", "A37", "A47")
    "B100" = @("B27:C27", "B37:C37", "B47:C47")
    "C100" = @("D27", "D37", "D47")
    "D100" = @("A28", "A38", "A48")
    "E100" = @("B28:C28", "B38:C38", "B48:C48")
    "F100" = @("D28", "D38", "D48")
}

foreach ($src in $applyMap.Keys) {
    $ws.Range($src).Copy()
    foreach ($dest in $applyMap[$src]) {
        $ws.Range($dest).PasteSpecial($xlPasteFormats)
    }
}

# Clear the scratch staging row.
$ws.Range("A100:H100").Clear()

# ---------------------------------------------------------------------------
# 3) Content edits: two task descriptions were reworded.
# ---------------------------------------------------------------------------
$ws.Range("B34").Value = "内容:初步编写android ui界面代码"
$ws.Range("B42").Value = "内容:找android和web端界面ui素材"

# ---------------------------------------------------------------------------
# 4) The saved selection moved from B43 to B34.
# ---------------------------------------------------------------------------
$ws.Range("B34").Select()
